$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 351-352; this pushes the existing rows
# 351-438 down to 353-440 (and everything else below shifts too),
# which already reproduces the "each record pair drops into the slot
# two rows down" pattern seen across the whole diff, including the
# brand new trailing rows 439/440 that end up identical to the old
# 437/438 ("Terminal La Palmera de La Serena" / Coliflor, 44491).
$ws.Range("A351:A352").EntireRow.Insert()

# Populate the newly inserted row 351 (Primera) with the new weekly
# record (same shape as its neighbours, new date + prices).
$ws.Cells.Item(351, 1).Value = 8
$ws.Cells.Item(351, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(351, 3).Value = "Coquimbo"
$ws.Cells.Item(351, 4).Value = 44508
$ws.Cells.Item(351, 5).Value = 4
$ws.Cells.Item(351, 6).Value = 100112008
$ws.Cells.Item(351, 7).Value = "Coliflor"
$ws.Cells.Item(351, 8).Value = "Sin especificar"
$ws.Cells.Item(351, 9).Value = "Primera"
$ws.Cells.Item(351, 10).Value = 2600
$ws.Cells.Item(351, 11).Value = 600
$ws.Cells.Item(351, 12).Value = 700
$ws.Cells.Item(351, 13).Value = 650
$ws.Cells.Item(351, 14).Value = "$/unidad"
$ws.Cells.Item(351, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(351, 16).Value = 650
$ws.Cells.Item(351, 17).Value = 1
$ws.Cells.Item(351, 18).Value = "Hortaliza"

# Populate the newly inserted row 352 (Segunda) with the matching record.
$ws.Cells.Item(352, 1).Value = 8
$ws.Cells.Item(352, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(352, 3).Value = "Coquimbo"
$ws.Cells.Item(352, 4).Value = 44508
$ws.Cells.Item(352, 5).Value = 4
$ws.Cells.Item(352, 6).Value = 100112008
$ws.Cells.Item(352, 7).Value = "Coliflor"
$ws.Cells.Item(352, 8).Value = "Sin especificar"
$ws.Cells.Item(352, 9).Value = "Segunda"
$ws.Cells.Item(352, 10).Value = 1600
$ws.Cells.Item(352, 11).Value = 500
$ws.Cells.Item(352, 12).Value = 550
$ws.Cells.Item(352, 13).Value = 525
$ws.Cells.Item(352, 14).Value = "$/unidad"
$ws.Cells.Item(352, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(352, 16).Value = 525
$ws.Cells.Item(352, 17).Value = 1
$ws.Cells.Item(352, 18).Value = "Hortaliza"
